# Update column F (dSF) values for the rows whose source data was re-pulled.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    3  = -4
    6  = -4
    8  = -6
    9  = -4
    14 = 4
    15 = 0
    18 = -3
    20 = -4
    26 = -3
    27 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
